$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Salary"
$ws.Range("G1").Value = "Owner"
$ws.Range("C1").Value = "Team"
$ws.Range("D1").Value = "Bye"
$ws.Range("F1").Value = "ESPN Projection"

$ws.Range("C1:G1").Font.Name = "Times New Roman"
$ws.Range("C1:G1").Font.Size = 12
